$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MAFMC ABC_ACL_catch")

# Columns that hold the Atlantic mackerel ABC/ACL values merged across
# row 11 (Rec) and row 12 (Comm).
$cols = @("B", "D", "F", "H", "J", "L", "N", "P", "R")

foreach ($col in $cols) {
    $topCell = $col + "11"
    $bottomCell = $col + "12"
    $mergedRange = $ws.Range($topCell + ":" + $bottomCell)

    # Grab the ABC/ACL value currently shown (and merged) on row 11.
    $abcValue = $ws.Range($topCell).Value2

    # Break the vertical merge between row 11 and row 12 for this column.
    $mergedRange.UnMerge()

    # Row 12 (Atlantic mackerel Comm) now shows the same ABC/ACL figure
    # as row 11 (Atlantic mackerel Rec) since they share one overall ABC.
    $ws.Range($bottomCell).Value = $abcValue

    # After unmerging, Excel drops the horizontal centering that only made
    # sense while the tall merged cell was centered; vertical centering is
    # kept.
    $ws.Range($topCell).HorizontalAlignment = 1
    $ws.Range($bottomCell).HorizontalAlignment = 1
    $ws.Range($topCell).VerticalAlignment = -4108
    $ws.Range($bottomCell).VerticalAlignment = -4108
}

# Reflect the selection the user ended up with after performing the edit.
$ws.Range("R11:R12").Select()
